$d = $word.ActiveDocument

# 1. Refresh the "last fetched" timestamp stamped in the footer.
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute(
    "2025-06-30 12:13Z / ", $true, $false, $false, $false, $false,
    $true, 1, $false, "2025-07-02 02:48Z / ", 2) | Out-Null

# 2. Add the regression-test character styles used to round-trip basic
#    inline HTML formatting (bold / italic / subscript / superscript /
#    underline) through Word.
$dpf = $d.Styles("DefaultParagraphFont")

$sb = $d.Styles.Add("b", 2)
$sb.BaseStyle = $dpf
$sb.Priority = 1
$sb.QuickStyle = $true
$sb.Font.Bold = $true

$si = $d.Styles.Add("i", 2)
$si.BaseStyle = $dpf
$si.Priority = 1
$si.QuickStyle = $true
$si.Font.Italic = $true

$ssub = $d.Styles.Add("sub", 2)
$ssub.BaseStyle = $dpf
$ssub.Priority = 1
$ssub.QuickStyle = $true
$ssub.Font.Subscript = $true

$ssup = $d.Styles.Add("sup", 2)
$ssup.BaseStyle = $dpf
$ssup.Priority = 1
$ssup.QuickStyle = $true
$ssup.Font.Superscript = $true

$su = $d.Styles.Add("u", 2)
$su.BaseStyle = $dpf
$su.Priority = 1
$su.QuickStyle = $true
$su.Font.Underline = 1
